$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Mike's hours (column D) for Week 8 (row 11)
$ws.Range("D11").Value = 28

# Update Patrick's hours (column F) for Week 8 (row 11)
$ws.Range("F11").Value = 20

# Add Mike's hours (column D) for Week 9 (row 12), previously empty
$ws.Range("D12").Value = 12

# Update the active cell selection to F12 to match the saved view state
$ws.Range("F12").Select()
